$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2 = @{ D=43.02000045776367; E=46.40999984741211; F=47.36999893188477; G=42.20000076293945; H=185748673 }
    3 = @{ D=47.09000015258789; E=49.88999938964844; F=50.20000076293945; G=46.43999862670898; H=185748673 }
    4 = @{ D=50.72000122070312; E=46.93000030517578; F=52.65000152587891; G=43.45999908447266; H=185748673 }
    5 = @{ D=50.09000015258789; E=50.08000183105469; F=52.77999877929688; G=47.72999954223633; H=185748673 }
    6 = @{ D=42.61999893188477; E=44.75; F=46.40000152587891; G=40.36000061035156; H=185748673 }
    7 = @{ D=47.52000045776367; E=51.66999816894531; F=51.9900016784668; G=46.68000030517578; H=185748673 }
    8 = @{ D=54.31999969482422; E=59.29000091552734; F=59.59000015258789; G=53.52999877929688; H=185748673 }
    9 = @{ D=59.25; E=60.47999954223633; F=62; G=57.54999923706055; H=185748673 }
    11 = @{ D=73.87000274658203; E=74.87000274658203; F=75.33999633789062; G=70.55000305175781; H=185748673 }
    12 = @{ D=76.9000015258789; E=80.41999816894531; F=80.44999694824219; G=75.52999877929688; H=185748673 }
    13 = @{ D=87.06999969482422; E=90.37999725341795; F=94.8000030517578; G=85.68000030517578; H=185748673 }
    14 = @{ D=91.88999938964844; E=84.66999816894531; F=92.94000244140624; G=82.0999984741211; H=185748673 }
    15 = @{ D=85.4000015258789; E=88.06999969482422; F=91.63999938964844; G=85.13999938964844; H=185748673 }
    16 = @{ D=89.2300033569336; E=102.1399993896484; F=102.9100036621094; G=89.0199966430664; H=185748673 }
    17 = @{ D=89.55000305175781; E=91.94000244140624; F=94.44000244140624; G=82.76000213623047; H=185748673 }
    18 = @{ D=93.51000213623048; E=101.6800003051758; F=104.1999969482422; G=92.9000015258789; H=185748673 }
    19 = @{ D=121.3199996948242; E=116.4400024414062; F=124.120002746582; G=113.9400024414062; H=185748673 }
    20 = @{ D=133.3500061035156; E=141.8099975585938; F=142.9400024414062; G=123.8499984741211; H=185748673 }
    21 = @{ D=136.3999938964844; E=141.0399932861328; F=143; G=133.8399963378906; H=185748673 }
    22 = @{ D=148.7700042724609; E=137.9299926757812; F=166.8699951171875; G=131.5099945068359; H=185748673 }
    23 = @{ D=153.6799926757812; E=180.9100036621093; F=181.9700012207031; G=147.9799957275391; H=185748673 }
    24 = @{ D=200.7400054931641; E=221.3000030517578; F=222.0299987792969; G=190.5; H=185748673 }
    25 = @{ D=215.5200042724609; E=227.5; F=246.6900024414062; G=211.1999969482422; H=185748673 }
    26 = @{ D=259.9599914550781; E=245.2100067138672; F=300.9100036621094; G=242.1300048828125; H=185748673 }
    27 = @{ D=248.9100036621093; E=254.3399963378907; F=259.1300048828125; G=225.0200042724609; H=185748673 }
    28 = @{ D=289.8200073242188; E=332.239990234375; F=335.010009765625; G=286.010009765625; H=185748673 }
    29 = @{ D=333.4500122070312; E=341; F=361.2200012207031; G=328.2699890136719; H=185748673 }
    30 = @{ D=311.8399963378906; E=312.3900146484375; F=324.5199890136719; G=275; H=185748673 }
    31 = @{ D=286.5499877929688; E=319.2000122070312; F=324.5; G=255.0200042724609; H=185748673 }
    32 = @{ D=364.1700134277344; E=346.0199890136719; F=391.1700134277344; G=341.6199951171875; H=185748673 }
    33 = @{ D=297.4500122070312; E=339.5400085449219; F=340.6099853515625; G=269.8800048828125; H=185748673 }
    34 = @{ D=353.6700134277344; E=363.760009765625; F=379.760009765625; G=347.9700012207031; H=185748673 }
    35 = @{ D=370.760009765625; E=454.9599914550781; F=468.0299987792969; G=365.4599914550781; H=185748673 }
    36 = @{ D=448.1400146484375; E=458.8900146484375; F=463.9700012207031; G=418.510009765625; H=185748673 }
    37 = @{ D=470.7000122070313; E=543.22998046875; F=564.780029296875; G=468; H=185748673 }
    38 = @{ D=534; E=573.72998046875; F=629.3800048828125; G=527.030029296875; H=185748673 }
    39 = @{ D=528.1699829101562; E=560.7999877929688; F=605.4500122070312; G=512.1199951171875; H=185748673 }
    40 = @{ D=556.0399780273438; E=519.5800170898438; F=585.6699829101562; G=474.5700073242188; H=185748673 }
    41 = @{ D=508.8299865722656; E=558.489990234375; F=574.6400146484375; G=508.8299865722656; H=185748673 }
    42 = @{ D=511.6900024414063; E=457.2799987792969; F=545.6500244140625; G=448.1099853515625; H=185748673 }
    43 = @{ D=462.4700012207031; E=463.9800109863281; F=521.0999755859375; G=437.6900024414063; H=185748673 }
    44 = @{ D=617.5499877929688; E=603.52001953125; F=636.6099853515625; G=593; H=185748673 }
}

foreach ($row in $data.Keys) {
    $rowdata = $data[$row]
    $ws.Range("D$row").Value2 = $rowdata.D
    $ws.Range("E$row").Value2 = $rowdata.E
    $ws.Range("F$row").Value2 = $rowdata.F
    $ws.Range("G$row").Value2 = $rowdata.G
    $ws.Range("H$row").Value2 = $rowdata.H
    $ws.Range("I$row").Value2 = "SNPS"
}

Write-Output "done"